$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'267.25"

$ws.Range("D3").Value = "'21.48"

$ws.Range("D4").Value = "'6.269"

$ws.Range("D5").Value = "'0.06167"

$ws.Range("D6").Value = "'3.570"

$ws.Range("D7").Value = "'6.584"

$ws.Range("D8").Value = "'1.379"

$ws.Range("D9").Value = "'0.8252"

$ws.Range("D10").Value = "'0.01347"

$ws.Range("D11").Value = "'0.1546"

$ws.Range("D12").Value = "'0.08212"

$ws.Range("D13").Value = "'0.03304"

$ws.Range("D14").Value = "'0.03183"

$ws.Range("D15").Value = "'0.09305"

$ws.Range("D16").Value = "'3.734"

$ws.Range("D17").Value = "'0.001623"

$ws.Range("D18").Value = "'0.04688"

$ws.Range("D19").Value = "'0.006304"

$ws.Range("D20").Value = "'0.005795"

$ws.Range("B21").Value = "UpBots"
$ws.Range("C21").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D21").Value = "'0.007489"
$ws.Range("E21").Value = "20UpBotsUBXTBestin24h"

$ws.Range("B22").Value = "BitKan"
$ws.Range("C22").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D22").Value = "'0.001068"
$ws.Range("E22").Value = "21BitKanKAN"

$ws.Range("B23").Value = "NitroEx"
$ws.Range("C23").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D23").Value = "'0.0001500"
$ws.Range("E23").Value = "22NitroExNTX"

$ws.Range("B24").Value = "LEO"
$ws.Range("C24").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D24").Value = "'3.725"
$ws.Range("E24").Value = "23LEOLEO"

$ws.Range("B25").Value = "BTSEToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D25").Value = "'2.319"
$ws.Range("E25").Value = "24BTSETokenBTSE"

$ws.Range("B26").Value = "BitpandaEcosystemToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D26").Value = "'0.3303"
$ws.Range("E26").Value = "25BitpandaEcosystemTokenBEST"

$ws.Range("B27").Value = "ProBitToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D27").Value = "'0.1242"
$ws.Range("E27").Value = "26ProBitTokenPROB"

$ws.Range("B28").Value = "AAXToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/LNePqkIhk+aaxtoken-aab"
$ws.Range("D28").Value = "'0.3999"
$ws.Range("E28").Value = "27AAXTokenAAB"

$ws.Range("D40").Value = "'0.04653"

$ws.Range("D41").Value = "'0.006983"

$ws.Range("D42").Value = "'0.004000"
$ws.Range("E42").Value = "41CEJICEJIWorstin24h"

$ws.Range("D43").Value = "'0.1137"

$ws.Range("D44").Value = "'0.01185"

$ws.Range("D45").Value = "'0.00005890"

$ws.Range("D46").Value = "'0.0009897"
$ws.Range("E46").Value = "45ACDXExchangeACXT"

$ws.Range("D47").Value = "'0.00000000750"

$ws.Range("D48").Value = "'0.7821"

$ws.Range("D49").Value = "'0.002442"

$ws.Range("D50").Value = "'0.00001900"

$ws.Range("D51").Value = "'0.01240"
